# Auto-generated edit script: updates numeric cell values per the commit diff
# for Sheets/Siren_Profits.xlsx (workbook sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5230.3335
$ws.Range("I51").Value = 4279.4
$ws.Range("K51").Value = 4279.4
$ws.Range("M51").Value = -3795.4
$ws.Range("H100").Value = 13949189
$ws.Range("I100").Value = 15662382
$ws.Range("J100").Value = 243648.5
$ws.Range("K100").Value = 15662382
$ws.Range("L100").Value = 243648.5
$ws.Range("M100").Value = -15661841
$ws.Range("N100").Value = -244730.5
$ws.Range("H106").Value = 4049.32
$ws.Range("I106").Value = 4456.1113
$ws.Range("K106").Value = 4456.1113
$ws.Range("M106").Value = -3825.1113
$ws.Range("H138").Value = 3810.018
$ws.Range("I138").Value = 707.8333
$ws.Range("J138").Value = 4675.744
$ws.Range("K138").Value = 2123.4999
$ws.Range("L138").Value = 14027.232
$ws.Range("M138").Value = 3016.5001
$ws.Range("N138").Value = -24307.232
$ws.Range("H140").Value = 80694
$ws.Range("J140").Value = 80694
$ws.Range("L140").Value = 80694
$ws.Range("N140").Value = -91054

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7759.95
$ws.Range("I2").Value = 8839.75
$ws.Range("K2").Value = 8839.75
$ws.Range("M2").Value = -8726.75
$ws.Range("H32").Value = 3300.8157
$ws.Range("I32").Value = 3440.8857
$ws.Range("J32").Value = 1666.6666
$ws.Range("K32").Value = 3440.8857
$ws.Range("L32").Value = 1666.6666
$ws.Range("M32").Value = -3153.8857
$ws.Range("N32").Value = -2240.6666
$ws.Range("H97").Value = 4517.3516
$ws.Range("I97").Value = 4664.4707
$ws.Range("K97").Value = 4664.4707
$ws.Range("M97").Value = -4168.4707
$ws.Range("H114").Value = 199000
$ws.Range("J114").Value = 199000
$ws.Range("L114").Value = 199000
$ws.Range("N114").Value = -207678
$ws.Range("H116").Value = 7759.95
$ws.Range("I116").Value = 8839.75
$ws.Range("K116").Value = 8839.75
$ws.Range("M116").Value = -6545.75
$ws.Range("H132").Value = 3856.838
$ws.Range("I132").Value = 2790.5417
$ws.Range("K132").Value = 8371.625100000001
$ws.Range("M132").Value = -5841.625100000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7759.95
$ws.Range("I3").Value = 8839.75
$ws.Range("K3").Value = 8839.75
$ws.Range("M3").Value = -8725.75
$ws.Range("H20").Value = 3254.077
$ws.Range("I20").Value = 1002.5
$ws.Range("K20").Value = 1002.5
$ws.Range("M20").Value = -755.5
$ws.Range("H105").Value = 3661.3635
$ws.Range("I105").Value = 2534.375
$ws.Range("K105").Value = 2534.375
$ws.Range("M105").Value = -787.375
$ws.Range("H110").Value = 267567.34
$ws.Range("J110").Value = 267567.34
$ws.Range("L110").Value = 267567.34
$ws.Range("N110").Value = -275747.34

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1101.4375
$ws.Range("I16").Value = 1005.0714
$ws.Range("K16").Value = 1005.0714
$ws.Range("M16").Value = -718.0714
$ws.Range("H68").Value = 99997.5
$ws.Range("J68").Value = 99997.5
$ws.Range("L68").Value = 99997.5
$ws.Range("N68").Value = -101495.5
$ws.Range("H71").Value = 99997.5
$ws.Range("J71").Value = 99997.5
$ws.Range("L71").Value = 299992.5
$ws.Range("N71").Value = -307480.5
$ws.Range("H113").Value = 1101.4375
$ws.Range("I113").Value = 1005.0714
$ws.Range("K113").Value = 1005.0714
$ws.Range("M113").Value = 1164.9286
$ws.Range("H135").Value = 63597.8
$ws.Range("J135").Value = 61997.5
$ws.Range("L135").Value = 61997.5
$ws.Range("N135").Value = -72137.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 477820.72
$ws.Range("I5").Value = 1474.4615
$ws.Range("J5").Value = 1251883.4
$ws.Range("K5").Value = 4423.3845
$ws.Range("L5").Value = 3755650.2
$ws.Range("M5").Value = -4311.3845
$ws.Range("N5").Value = -3755874.2
$ws.Range("H68").Value = 55570556
$ws.Range("I68").Value = 2497.5
$ws.Range("K68").Value = 7492.5
$ws.Range("M68").Value = -6681.5
$ws.Range("H71").Value = 55570556
$ws.Range("I71").Value = 2497.5
$ws.Range("K71").Value = 22477.5
$ws.Range("M71").Value = -18421.5
$ws.Range("H98").Value = 740.8125
$ws.Range("I98").Value = 947.1429000000001
$ws.Range("J98").Value = 580.3333
$ws.Range("K98").Value = 2841.4287
$ws.Range("L98").Value = 1740.9999
$ws.Range("M98").Value = -1343.4287
$ws.Range("N98").Value = -4736.9999
$ws.Range("H129").Value = 3889.7
$ws.Range("J129").Value = 6999.6
$ws.Range("L129").Value = 20998.8
$ws.Range("N129").Value = -30998.8
$ws.Range("H135").Value = 477820.72
$ws.Range("I135").Value = 1474.4615
$ws.Range("J135").Value = 1251883.4
$ws.Range("K135").Value = 13270.1535
$ws.Range("L135").Value = 11266950.6
$ws.Range("M135").Value = -10735.1535
$ws.Range("N135").Value = -11272020.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 9008000
$ws.Range("I18").Value = 27000000
$ws.Range("J18").Value = 12000
$ws.Range("K18").Value = 27000000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = -26999707
$ws.Range("N18").Value = -12586
$ws.Range("H70").Value = 9299.5
$ws.Range("I70").Value = 4933
$ws.Range("K70").Value = 4933
$ws.Range("M70").Value = -4663
$ws.Range("H73").Value = 9299.5
$ws.Range("I73").Value = 4933
$ws.Range("K73").Value = 4933
$ws.Range("M73").Value = -3997
$ws.Range("H103").Value = 46151
$ws.Range("J103").Value = 46151
$ws.Range("L103").Value = 46151
$ws.Range("N103").Value = -48495
$ws.Range("H122").Value = 23699.5
$ws.Range("J122").Value = 20673.8
$ws.Range("L122").Value = 62021.39999999999
$ws.Range("N122").Value = -66921.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2339.3044
$ws.Range("I22").Value = 2614.5715
$ws.Range("K22").Value = 2614.5715
$ws.Range("M22").Value = -2319.5715
$ws.Range("H23").Value = 6666.6665
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 9000
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = -1770
$ws.Range("N23").Value = -9460
$ws.Range("H27").Value = 2339.3044
$ws.Range("I27").Value = 2614.5715
$ws.Range("K27").Value = 2614.5715
$ws.Range("M27").Value = -2507.5715
$ws.Range("H136").Value = 10404.866
$ws.Range("I136").Value = 8725.714
$ws.Range("K136").Value = 26177.142
$ws.Range("M136").Value = -23627.142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 20279.857
$ws.Range("I122").Value = 3341.842
$ws.Range("J122").Value = 40393.75
$ws.Range("K122").Value = 10025.526
$ws.Range("L122").Value = 121181.25
$ws.Range("M122").Value = -7575.526
$ws.Range("N122").Value = -126081.25
$ws.Range("H126").Value = 25219.21
$ws.Range("I126").Value = 35781.668
$ws.Range("K126").Value = 107345.004
$ws.Range("M126").Value = -104875.004

